$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6440
$ws1.Range("F3").Value = 115
$ws1.Range("F5").Value = 385
$ws1.Range("F7").Value = 5
$ws1.Range("F13").Value = 372
$ws1.Range("F14").Value = 942
$ws1.Range("F15").Value = 3148
$ws1.Range("F17").Value = 191
$ws1.Range("F18").Value = 1820

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6440
$ws4.Range("F3").Value = 115
$ws4.Range("F5").Value = 385
$ws4.Range("F7").Value = 5
$ws4.Range("F14").Value = 372
$ws4.Range("F15").Value = 942
$ws4.Range("F16").Value = 3148
$ws4.Range("F18").Value = 191
$ws4.Range("F19").Value = 1820
